$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value2 = "66.540.52"
$ws.Range("E2").Value2 = "  +0.82%  "
$ws.Range("D3").Value2 = "3.599.88"
$ws.Range("E3").Value2 = "  +1.59%  "
$ws.Range("E4").Value2 = "  -0.16%  "
$ws.Range("D5").Formula = "'609.52"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value2 = "  +0.65%  "
$ws.Range("D6").Formula = "'148.60"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value2 = "  +2.83%  "
$ws.Range("E7").Value2 = "  +0.05%  "
$ws.Range("E8").Value2 = "  -0.71%  "
$ws.Range("E9").Value2 = "  +1.96%  "
$ws.Range("E10").Value2 = "  +0.07%  "
$ws.Range("E11").Value2 = "  +1.06%  "
$ws.Range("D12").Value2 = "4.210.00"
$ws.Range("E12").Value2 = "  +1.53%  "
$ws.Range("E13").Value2 = "  +1.33%  "
$ws.Range("D14").Formula = "'29.82"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value2 = "  -0.56%  "
$ws.Range("D15").Value2 = "3.606.53"
$ws.Range("E15").Value2 = "  +1.61%  "
$ws.Range("D16").Value2 = "66.654.54"
$ws.Range("E16").Value2 = "  +0.63%  "
$ws.Range("E17").Value2 = "  +0.81%  "
$ws.Range("D18").Formula = "'11.51"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value2 = "  +1.80%  "
$ws.Range("E19").Value2 = "  +3.48%  "
$ws.Range("D20").Formula = "'15.11"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value2 = "  +1.87%  "
$ws.Range("D21").Formula = "'428.18"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value2 = "  -0.44%  "
$ws.Range("E22").Value2 = "  +1.70%  "
$ws.Range("E23").Value2 = "  -0.41%  "
$ws.Range("D24").Value2 = "3.743.91"
$ws.Range("E24").Value2 = "  +1.51%  "
$ws.Range("E25").Value2 = "  +0.02%  "
$ws.Range("D26").Formula = "'0.0000123"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value2 = "  +4.61%  "
$ws.Range("D27").Formula = "'8.34"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value2 = "  +5.12%  "
$ws.Range("D28").Formula = "'9.48"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value2 = "  +4.64%  "
$ws.Range("E29").Value2 = "  +0.12%  "
$ws.Range("D30").Formula = "'1.00"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value2 = "  -0.02%  "
$ws.Range("D31").Formula = "'1.48"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value2 = "  +0.92%  "
$ws.Range("D32").Value2 = "3.596.49"
$ws.Range("E32").Value2 = "  +1.53%  "
$ws.Range("D33").Formula = "'0.158"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value2 = "  +2.66%  "
$ws.Range("D34").Formula = "'25.49"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value2 = "  +0.18%  "
$ws.Range("D35").Formula = "'7.87"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value2 = "  +0.11%  "
$ws.Range("E36").Value2 = "  -0.01%  "
$ws.Range("D37").Formula = "'5.66"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value2 = "  +1.17%  "
$ws.Range("E38").Value2 = "  -1.87%  "
$ws.Range("D39").Formula = "'177.23"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value2 = "  +1.37%  "
$ws.Range("E40").Value2 = "  +0.87%  "
$ws.Range("E41").Value2 = "  +1.18%  "
$ws.Range("D42").Formula = "'0.900"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value2 = "  +1.07%  "
$ws.Range("E43").Value2 = "  -0.42%  "
$ws.Range("E44").Value2 = "  +9.84%  "
$ws.Range("D45").Formula = "'0.999"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value2 = "  -0.09%  "
$ws.Range("D46").Formula = "'25.06"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value2 = "  -2.04%  "
$ws.Range("D47").Formula = "'1.18"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value2 = "  -2.13%  "
$ws.Range("D48").Formula = "'24.07"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value2 = "  +3.09%  "
$ws.Range("E49").Value2 = "  +1.23%  "
$ws.Range("E50").Value2 = "  +1.55%  "
$ws.Range("D51").Value2 = "2.430.35"
$ws.Range("E51").Value2 = "  +5.69%  "
